# Updated output files to reflect 9 holdouts (was 11): McAdams (UT-04) and
# McBath (GA-06) are no longer counted among the impeachment "NO" votes.
#
# - top_trump_dists: flip their for_impeachment flag to YES and stamp the
#   vote date (2019-10-04) onto the row.
# - full_list_of_nos: remove their rows outright (the list shrinks from 11
#   to 9 rows) and shift the remaining holdouts up.
# - the various groupings_* summary tabs: rebalance the NO/YES counts that
#   moved because of the flip.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. prezresults2016 — R/NO count drops by 2, R/YES count gains 2
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("prezresults2016")
$ws.Range("C3").Value = 8
$ws.Range("C5").Value = 23

# ---------------------------------------------------------------------
# 2. gdp_vs_nationalavg
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("gdp_vs_nationalavg")
$ws.Range("C2").Value = 2
$ws.Range("C4").Value = 131

# ---------------------------------------------------------------------
# 3. college_vs_nationalavg
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("college_vs_nationalavg")
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 8
$ws.Range("C4").Value = 134
$ws.Range("C5").Value = 92

# ---------------------------------------------------------------------
# 4. nonwhite_vs_nationalavg
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("nonwhite_vs_nationalavg")
$ws.Range("C2").Value = 3
$ws.Range("C3").Value = 6
$ws.Range("C4").Value = 141
$ws.Range("C5").Value = 85

# ---------------------------------------------------------------------
# 5. rural_morethanfifth
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("rural_morethanfifth")
$ws.Range("C3").Value = 4
$ws.Range("C5").Value = 191

# ---------------------------------------------------------------------
# 6. margin_5_or_less
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("margin_5_or_less")
$ws.Range("C2").Value = 6
$ws.Range("C4").Value = 20

# ---------------------------------------------------------------------
# 7. top_trump_dists — flip McAdams (row 11) and McBath (row 26) to YES
#    and record the date they came out for impeachment.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("top_trump_dists")

foreach ($r in @(11, 26)) {
    $ws.Range("A$r").Value = "YES"

    $ws.Range("H$r").Value = 43742
    $ws.Range("H$r").NumberFormat = "yyyy-mm-dd"

    $ws.Range("J$r").Value = 43742
    $ws.Range("J$r").NumberFormat = "yyyy-mm-dd"

    $ws.Range("K$r").Value = 10
    $ws.Range("L$r").Value = 2019
}

# ---------------------------------------------------------------------
# 8. full_list_of_nos — drop McAdams (row 8) and McBath (row 9); the
#    remaining holdouts (Peterson, Torres Small, Van Drew, ...) shift up
#    two rows so the table shrinks from 12 to 10 rows.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("full_list_of_nos")
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()
